$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Ocupacao"
$ws.Range("B1").Value = "Pesos"
$ws.Range("C1").Value = "Matriz_Covariancia"
$ws.Range("D1").Value = "Erro_Estimacao_Amplitude"
$ws.Range("E1").Value = "Media_Erro_Estimacao"
$ws.Range("F1").Value = "Desvio_Padrao_Erro_Estimacao"

$headerRange = $ws.Range("A1:F1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$ws.Range("A2").Value = 20
$ws.Range("B2").Value = "[-0.14339454108469085, 0.11960066703238818, -0.9949803054570084, 2.1723204512071397, -1.4340851491821278, 0.670651331924537, -0.3901124544402377]"
$ws.Range("C2").Value = "[[ 5.05233430e+02  3.59521213e+02  1.41835422e+02  3.72691888e+01`n   7.31522850e+00  9.03398774e-01 -4.08527618e-02]`n [ 3.59521213e+02  4.98091789e+02  3.55207374e+02  1.41522106e+02`n   3.84938636e+01  7.45866855e+00 -1.37401295e-01]`n [ 1.41835422e+02  3.55207374e+02  5.02425347e+02  3.62997103e+02`n   1.44129886e+02  3.69871686e+01  5.32898808e+00]`n [ 3.72691888e+01  1.41522106e+02  3.62997103e+02  5.13320673e+02`n   3.65379898e+02  1.43639191e+02  3.77415259e+01]`n [ 7.31522850e+00  3.84938636e+01  1.44129886e+02  3.65379898e+02`n   5.09426324e+02  3.64256507e+02  1.45017509e+02]`n [ 9.03398774e-01  7.45866855e+00  3.69871686e+01  1.43639191e+02`n   3.64256507e+02  5.10482666e+02  3.62599296e+02]`n [-4.08527618e-02 -1.37401295e-01  5.32898808e+00  3.77415259e+01`n   1.45017509e+02  3.62599296e+02  5.03142586e+02]]"
$ws.Range("D2").Value = "[10.963804422385438, -6.419417684714378, 8.591147980621551, -6.909173094399636, 6.654634230543489, 0.262950001603334, -2.3177950756164174, 3.7949114617889412, -1.6758122559240234, -0.029768978717908068, 4.014176256208368, -0.9531122546353448, 7.093454556998203, 53.809865849359056, 26.653957798671954, 51.39348811293895, 32.65744510146598, 3.986092007957211, -11.813679641368909, 86.96642169430693, 61.99121724958726, 61.417224694271965, 11.369343236079029, 9.447396359588822, 8.61314094119848, 25.32722868710443, 30.528222848148417, 20.841381471268022, 10.210560744143208, -1.3980847546838024, 3.0695674465896143, -1.7747976561966217, 2.849855725199466, -2.6821364788362145, -1.948949774707485, 4.445174069940222, -5.658129696904189, 8.48745064910769, -10.124341229837459, 7.4579893792946, -2.8412299537658567, 28.456181601917823, 2.1479150205755317, 33.386754324713294, -1.9809703549019524, 26.877139382619774, 16.114128071228883, 51.87593141529892, 45.102299909789245, 91.71878799623234, 11.161513049923286, 133.0313549697823, 27.40504829822924, 27.94219753824035, 19.413164710834977, 98.85967729254322, 58.841096079212846, 23.334583072972336, 27.11674883890105, 7.231341266378159, 10.258779876908426, -1.6383544817196705, 5.483560648565518, -6.949225144983137, 8.812122051453294, 4.144210959939853, -2.9228163298876346, 14.054664552093346, -1.7249108508642053, -5.383841276492426, 4.067849440687723, 13.507613145109124, -6.239086167997154, 11.196909574010197, 20.29758595259385, 45.83399742016664, 7.228271877745925, 77.14311603227831, 8.782752268058772, -56.23847402274225, 51.195354044877654, 36.25364199684658, 25.42629073635089, 4.614173131366974, 4.788709500316826, -6.131101580910678, 4.216078748321351, 1.673180712322694, -3.7319317192988604, 0.2869401106253251, 1.2170301550327969, -3.439595722250816, 5.236145183917018, -2.726723600739994, 0.7471836172418019, -0.35633697779790724, -0.8348389050712577, 11.563735978989609, 8.43116536373828, 10.21139124141288, 20.340191683038455, -8.980877997315247, -1.2877298031584434, 6.787924533090226, 16.11081108233159, 8.223501755668778, 0.25100628214789794, 9.03151345742128, 6.856848348827693, 9.189043765223737, 20.412357633740015, 9.87368589093272, -22.24069357626605, 62.370558334509944, 36.231491249304966, 47.77771314552982, 41.07970174719967, -3.762467673612946, -13.899823846511211, 55.58210719170114, 35.79237150453775, 25.734854746506308, 1.7283701142715777, 4.801524699697631, -3.0125563251688154, -0.3329854163150259, 4.434475812580423, -1.4439573561749208, 8.703206816959497, 20.516922578258086, 5.906729985283036, 51.85206489444627, 36.45584835106372, 44.662513541873935, 55.65374284175736, 68.51504856694851, 5.104824859852869, 97.5284459181121, 73.5537785569966, 50.73627304671919, 2.6719362137107785, 62.085451968506376, 10.530828941196406, 33.03751743589956, 12.513477432249408, 20.226193963860318, 26.211197816677437, 19.233770498037586, 85.33300660636833, 18.2933030828992, 45.451180474899104, 0.671116113184155, 52.57744695599746, 22.740238178649587, 49.63775479007014, 47.3898186327428, 36.65153764112745, 37.979384947799694, 13.109728578987514, 38.257537678349536, 1.5953853616442046, -14.660526752895265, 66.48422293440566, 34.95657838786495, 34.43649749190274, 61.4841918448967, 68.36359935155949, 12.720756360771759, 134.96793979080866, 48.82213073925834, -50.428727803547105, 108.70884680164113, 38.73886873340189, 34.06583296633009, 24.597870065515718, 23.642155689697226, 11.227681205848034, 14.480429579338235, -2.3730794373736224, 5.53713319237124, 9.62323198415842, -3.9005252955596346, 2.8127587802556278, -6.89976918044721, 11.784067461410968, -3.027320327842644, 6.339343088582748, 23.758778625798474, 2.725419172498185, 28.42679655385639, -6.642065378407965, 4.5081557348771355, 6.51810031100843, 20.862940886987168, 14.352708767816798, 7.957536399803486, 1.763030785883009, -3.611079014369215, 9.275054100954812, 0.8032161810131715, 10.864302880382471, 10.20677393076912, 49.111447993594545, 37.3803474218615, 3.164461246587912, 81.80963387984023, 26.01106671046489, -54.786040580160325, 60.841619613837494, 30.428102622828405, 29.214102009818433, 14.372792887768673, 1.9236117359083167, 22.51002096123455, -5.6537346367639945, -0.8064812234167817, 5.1214371384684485, 7.672238602637837, 11.560652172435788, 2.6746901454245195, 0.16811660915466708, 5.46471290261834, 1.6858156400159823, 0.4288935179988598, 4.349143124932903, -4.8497480847844, 5.107480622743742, 3.827691185073015, 3.864022930698831, 1.7303702624289046, 2.7816986253773437, 1.7622950886711932, 0.9154552685812831, -2.6956466942265136, 1.9077778718086607, 2.8922364505300115, 1.9971651505634354, 1.5894752422200473, 14.12772872945139, 0.14220393286140975, 24.063877286073797, -7.6034401507461205, 13.67269665423968, 0.5119917134309315, 12.765007901164953, 8.305583869982264, 47.455628712492356, 22.69594376558142, 35.717158013867966, 26.449702100859184, 0.32679382525449796, -14.342018235891063, 44.26309187426695, 27.269380666876728, 18.92661740337368, 6.5904982134509, -1.942846485640219, 5.2978896246523846, 26.592667837472824, 18.122784342488572, 26.61161192944742, 22.87348069014812, 29.50555635971878, 13.991792584531545, 64.70106156297116, 62.4952942792778, -5.7429396318964905, 19.95763544769794, 29.646249659004283, 41.09474068066718, 24.92590090935548, 9.090439126403512, 9.969289584845132, 4.111122531819541, 21.643070320056726, -4.564298110319754, 26.159065590953563, 12.22646466425677, 20.641584570725982, 39.19203245692782, 13.01443498809456, 14.332250716611682, 3.5725477942432264, 44.16978156079248, 39.10987424348361, 17.50803219753817, 44.82933776563355, 9.381375446574694, -37.81055168109461, 36.069606458521115, 7.480051561439793, 22.674391953409554, -3.4254575465872357, 6.7797013998456555, -1.339268288830581, -2.030867217461452, 7.941542608092689, -13.50057476111568, 9.925101253132055, 1.7699695493157392, 18.632801274711113, 2.6255934973794197, 24.48206139124677, 2.889902768331595, -22.73610868080233, 25.546249644461184, 9.944702647763862, 17.93404062627791, 7.322012451263397, 19.831283322910842, 51.53835767337931, 14.881030571421562, 66.17775624563444, -8.726691356930354, 26.29572142072027, 16.01075050272229, 58.66326255545397, 36.521781330232265, 27.306187925251752, 27.779316302747702, 10.223999069751343, 29.062581476765708, -2.2058410822553967, -3.4011319345563606, 20.592567103897732, 17.367627671341747, 14.27872314897354, 9.381297993811822, -1.8555961246797155, 1.212212324974744, -0.18376428570115602, 1.2522895030552288, -3.029845853588345, 4.618541348781362, 18.54151005321054, 38.92740179276902, 25.3868076120717, 55.67275137323004, 14.258622151404339, -15.611362023552203, 34.89151751140287, 56.451292542961895, 57.01302809495756, 14.311724909377462, 2.1675148696117077, 3.0636097837153784, 29.121419313641276, 14.374953197742647, 31.9358115547629, 32.66001453974671, 9.403567526705494, 61.150499280172184, 4.701128658043899, -48.97986299924297, 45.82242525035449, 27.09774706853283, 15.057403491709286, 10.963486276871636, -3.841606962764413, 3.300642101360493, -0.9812985578081236, -1.9386700837531026, 1.5062811466568613, 4.318209413802614, 9.489993676296358, 14.264017129042365, 26.205408731762304, 21.18373190421681, 17.410472521219603, -19.85809410459829, 27.102376255896814, 25.86063759110166, 18.609492222317357, 43.451259461102964, 20.72249610425556, 16.235047370378147, 47.02050550851601, 8.342255917540342, -30.76337992519582, 41.263987126914564, 21.70091225107124, 42.42181796268072, 28.561534569973446, 24.858494627954002, 30.926734941137198, 0.8130217980017278, -18.91598511788269, 48.08167568506626, 51.714072243114124, 32.04211352896489, 36.47068295764774, -6.188925082967806, 4.691521740657464, 35.31951688104043, 45.03630853570855, 42.87085869824372, 16.397498805789027, 4.946714936871729, -4.514589879289554, 33.30150219693623, 14.612217262830415, 18.425109179406196, 16.854963566938622, 7.326099128849769, 17.421160093010258, 18.173884575885047, -2.783647168005004, -11.446108290072482, 24.228234091954267, 11.133028841824366, 31.859562417846213, 8.329491705456217, 31.102755417530652, -2.7124313338280004, 4.340950314098947, 4.115427660831929, 21.111264105042434, 16.29684613313592, 8.594879358757913, 4.7723524280366, 1.0960914346803592, 0.6307536665499364, 14.563679985686496, 7.6856163317956785, 14.166667443780295, 14.946205631332814, 0.18711464265571465, -8.996668696445509, 15.312197031972842, 30.02289785122899, 11.120740665612534, 27.198743940369926, 14.069593592029403, 11.399682516770994, 1.5895364751222898, 20.313552641793386, 2.3861631733468935, 19.119742723194296, 12.78668072308835, 4.512044134011193, -1.0509962779359379, 4.867827822838006, 0.7977792968374384, 16.1083911145095, 8.086329875656576, 10.581826120831554, 4.459924834969229, -11.374138399571741, 20.060896043754255, 3.8793189344019323, 7.476330353704768, 7.7675561705809075, -6.311219759104397, 8.177891510411236, -8.700401307650782, 7.185932653436935, -6.967029832405512, 7.19499932160935, 16.816910188496138, 9.797774531589688, 30.031230551656606, -3.2057851528208836, 4.1232926191771355, -1.880598105784996, 27.292589032259293, 12.569544896069628, 7.106430360878286, 5.916512340105086, 4.20236578205249, 62.33133314943409, 24.876747775558243, 98.6698572986256, -4.575486128171008, 32.57203156758396, 50.96485566408073, 64.72139254237442, 45.65309471828126, 55.57142083096507, 71.37510409999456, 32.47352476653593, 85.1955438011851, -34.30211604321703, 21.512767226466593, 48.95341035886932, 69.40768813059566, 70.33743675811593, 44.17444935613969, 1.5516421440417787, -12.731200262664107, 39.33586263105635, 36.467673635472856, 19.927290880127302, 10.10061885720237, 15.967142024803163, -0.9176835209061842, 4.785454586738169, 24.493384706251938, 11.39052641833576, 49.270464824631276, -5.852269263876553, 7.108178018505612, 4.522845043055673, 26.705186805197673, 11.719271188538729, 13.922920993416463, -0.3619554078176276, 3.794967038919492, -2.09594482148852, 0.038212215764959756, 0.7543538390741511, -3.721932903575941, 2.259223920031712, -1.3657481808214091, 5.311392578220338, -0.9563789954061519, -3.383591915898333, -1.4092593992358826, 2.183050112005874, 0.39775715801125244, -0.364427456677652, -0.23373924457681666, 0.9507203273638245, -3.4917064030092178, 2.056750722164607, 2.127247068810753, -1.7500820142228632, 3.061090970143044, -2.560772896933708, -0.5498679837659437, -0.5018550718716269, 0.7194361617006007, 0.18572509208188853, 3.9256757833546114, -5.268733106220703, 3.8296236803907524, -3.186283791111265, 47.92493611795526, 14.390824814891161, 57.18673126386305, 2.36198708940503, 3.4000219800295213, 2.278374976924689, 52.51206325430047, 36.881433036051384, 37.54310961745041, 11.686942417478768, 39.05125208148441, 10.737190277553324, 5.22518932536504, 18.950382225796545, 15.5780513906587, 12.524108107721403, 29.986620808347496, 15.66656466450523, 19.51906952400509, 3.5069212191571992, -0.4757855675860867, 16.339303066191484, 10.965097279066942, 11.49347498978327, -12.716992329905972, 9.969963353292963, 2.9214365857395337, 27.24583566369341, 25.90330038063462, 10.959796244776548, 52.867878039638484, 10.901475625110976, -41.732740084148475, 37.60899443815548, 24.54104491332751, 12.76127706542887, 14.27522726239414, -3.2098397333417488, 3.129460251496576, 0.1895534365347641, -0.8376076255928464, -2.462461138847518, 4.369994230515522, -5.459570012147273, 5.71693745132891, -6.259833151006488, 7.475119153676498, -6.7403976673748405, 4.851161325616529, -2.257730570446203, -2.479309512502845, 16.27965170871085, 2.901442530659949, 12.182093116463511, 12.525761171412885, -2.897325082659222, -6.955689966885845, 17.306880696437773, 5.649888189046847, 9.257683080878788, 2.1893780405379877, 1.5273751315339368, -1.310698509917214, -0.25141499417080093, 0.6502444259078368, 3.7870988453115992, 3.5652523532208678, -0.5729916594838279, 0.35714584649748815, 0.4541811290783624, 7.511882672059571, -4.849462957315222, 19.344810784318987, 3.423564401903633, 74.98481927818568, 89.14330945590538, 26.781905194418428, 132.5068212423966, 39.74379119373174, -109.88155380612822, 118.2481702285384, 67.55819356372366, 50.5207080515278, 33.58246450166527, 21.080671105101707, 16.56595412364965, 27.675908526273606, -3.2446853356833856, -3.9077179371201614, 30.818362045655245, 18.059155420960607, 24.838896306893314, 27.130733500179687, 1.1718980105541021, 36.33779899500191, -3.959806388455533, -14.981409516980822, 83.04345760058146, 34.829640340708565, 85.27128054006693, -6.621536296108122, 3.129578581665214, 16.459954889265614, 63.99621114141256, 40.095735825406784, 27.35293134073448, 11.434376719240841, 0.1664410395961866, 8.752137671236818, 4.660636378766872, 1.6752812508682133, 45.5622847967967, 18.168286504204502, 71.85930317786386, 1.871193364827171, 22.581979688156537, 6.277564022456815, 57.17968529881715, 38.98876576160699, 32.92460729893822, 54.49298216248907, 43.84990721853758, 126.04935821267863, 41.664081563767866, 142.08667202253656, 58.36616538913404, -39.75914258572652, 88.40945079905359, 148.7813657104155, 85.74457806944847, 28.87125846121046, 48.14151308502555, 40.712025156362586, -10.085872142719888, 38.049418492454976, 17.049044310239246, 24.757711725775835, 2.5569053420625387, 21.89464402795654, 1.2175054263138243, 14.746715891822127, -2.865212998948401, -0.7737702481971208, 16.456341858326322, 7.959308532120553, 12.957680245401498, 2.9725822671041495, 4.5218664339816605, -4.935067125612017, 2.2131240325296915, 3.6245355603642193, 1.3804583531960533, 24.496380141047677, 6.786818031962795, 30.11055239010064, -16.295869452957117, 11.166858697453705, 3.110660150521165, 27.407502081054353, 10.59282077863541, 7.393185437824538, 5.507318359147138, 8.042877599767003, 18.85125663903551, 28.370644071068178, 58.73552768045741, 86.97037306782097, 64.30077955005399, 119.04395421576336, 88.01018237952611, 79.91735362287329, 62.02043045737385, 149.8166989125922, 103.6966646773953, 77.70077482884722, -16.774142332553993, 90.0001007452624, 43.620969733449726, 31.09011966476024, 7.050925866751468, 10.27279809445509, 1.0341878690208806, 1.8413396425912634, -7.808699023649386, 13.48723037287959, -8.746408136188661, 9.146486750798845, -5.8782991632518655, 2.245984260020897, 3.3726576826912176, 0.748718881086484, 5.890697553204161, 10.411016195297712, 2.7601259770211364, 1.0762397988680803, -2.8373011882814874, 9.394875088170483, 2.2640258039997105, 8.255417114722412, -5.51876443012742, 7.480362062042602, -7.266961622773849, 9.186070732971494, -11.036386763269944, 8.2739376392174, -3.426411656753448, 1.6347063877521584, 30.558507958212818, 6.056996331385591, 32.50688369480063, 8.045485297318955, 1.2061224922382792, -8.24653392335716, 38.170339729711806, 14.069973250111007, 14.53515590440337, 6.0167530188606015, 43.206175742087005, 34.41234056600142, 37.561419003521216, 51.38165981121489, 5.871981738685577, -43.92262693537896, 72.50348770547032, 33.54159230109732, 30.8326839443879, 8.05989147559968, 2.60969934610743, 26.68877652808804, 11.135594687586135, 29.24222332181554, 14.379621426103682, -2.6603064577872786, -3.076542350633499, 25.9568801029108, 24.083485377784406, 37.40222329671735, 42.09664524782929, 37.04475182550744, 24.141147659240126, 22.415090713392743, 41.40443363472719, 70.8848647443081, 96.37217794447824, 63.81005469959287, -28.091418952441444, 14.713500450022973, 58.348359696537436, 46.376366816073, 29.562994724754226, 15.704569614645688, 23.375071559066097, 12.081660348973386, -1.4230330734703145, 37.517175244411106, 15.805617261165844, 0.5504919767299157, 24.602959543871734, 12.465540852789625, 20.902218287111566, 6.9160621698991935, 2.5991523980431452, -2.0082119480369442, 10.321064592332476, 15.338806685650274, 18.45862196628877, 6.821344339056253, 21.334664381908354, 4.677385872848092, -17.422818727325687, 26.740671598077153, 13.122260728249438, 17.757261979892956, -3.3172038586451422, 13.59330269902419, -12.909784764611548, 16.25724252218062, 3.06806234192549, 2.3215871668116383, -0.18542607949998002, 3.6354894413559293, 54.10550566812083, 61.991470144663595, 32.26556657417687, 122.92062176129146, 30.369965150339226, -98.83442911983367, 101.0378553583163, 70.12373319138246, 34.04583775842246, 12.128679863370477, 12.613450152062498, 6.82980414403246, 8.286031306165425, 6.781535908516424, 34.09516254675364, 40.86478719022034, 12.092013018864272, 65.88195537076346, 4.683202260255598, -8.060493474319646, 56.53340948773249, 62.99678291216307, 41.81712043979095, 28.936226362819582, 78.89283808861954, 83.70874131381474, 44.66326407890307, 92.38777683477234, 44.608673013745715, -22.58546857944628, 83.13056387805715, 37.65904820812647, 45.15834236878267, 33.77396072796406, 19.96137695138697, 22.030273680579263, 11.015547734156767, 7.026780337493264, 4.5518706912210325, 27.537472528181894, 12.237832532110492, 20.70867170300941, -1.510347341614737, 10.778225686203191, 1.7083804020953186, -0.01869015171678895, 0.12935111754559148, 1.108195122777898, 0.8759070893103593, -2.9213331837764165, 1.800146319868821, 6.657890902693673, 21.535511679424964, 1.3918333918494241, 29.310315296108456, 11.592626745703896, -29.987926738694284, 19.720948415139013, 10.075924427262445, 16.95577453610666, -2.961307190545028, 1.603887634048637, 6.172151771516482, -7.571419001200001, 8.643923712403197, 11.64885759750716, 31.476473262337457, 1.5588952798263165, 94.74111457123448, 67.21170103657474, 23.514627414497255, 133.17948651711828, 36.46734109297586, 10.516684614782037, 76.11962348451203, 70.10006646215412, 83.75746776103622, 59.6371345944251, 27.689873853910157, 32.132216432460474, 51.51389906889004, 24.880556919001535, 29.072953722163394, 26.228496994206125, 25.38871828986724, 40.234563980622816, 17.601731490070513, 6.853632712315957, 15.552379085911378, 15.086737955318773, -3.7949097537349252, 5.807741969498272, 2.6092444141494253, -2.5236650273872767, -0.313652376709598, 1.6217185830434797, -2.998734187671086, 5.63992475272879, -1.0724448868000873, -1.3151340588921898, -1.0362558361752228, 0.8267726482017168, 2.4346809754984244, 2.9284784271776996, -4.476455721668018, 5.14196410041807, -5.848238782953469, 20.2472427193008, 7.503363490580735, 43.4332945674136, 41.67058425375028, 35.56174728486343, 35.17944140867541, 37.784325368749876, 0.3188567949374175, 46.00617063276987, 69.40439002678448, 46.62247943554402, 41.26052552633463, 16.58303243169192, -22.684945852814934, 42.2955477129314, 12.33191029785181, 19.75162958413654, -3.4470384410390817, 8.294147140903513, 7.747935142231952, 6.552130642489152, 6.432672373121916, 1.8346242796124343, 2.799378273187912, 1.0551115370000375, -3.7256796723599277, 9.600394337437052, -0.980788974568144, 2.75530549353078, 5.287242036656306, -2.334154124797226, 5.881121098230103, 46.50222926692312, 12.958987865231421, 69.84134752324636, -12.469918320086045, 14.224445483350564, 8.164558324811061, 84.3475879143082, 38.3893166649647, 95.65470120416634, 2.1832530007062303, 73.86242486016297, -16.181083886300577, 39.32713504495533, 27.85352959879806, 48.565184725666654, 29.623747033926584, 33.08370384301778, 8.10345594495031, 12.679402558588482, 7.193041917390195, 4.351938047528677, 7.168231722846757, 17.117574596608396, 7.216575781077054, 2.831210427315554, 6.731950130052189, 6.800909117516573, -0.6487533213572941, 8.199956212972227, -5.000564494756847, 11.446129209292348, -5.741231434894834, 27.50083823117454, -0.03256431681068861, 16.56513004497762, -0.19537019931522526, 2.603253895291843, 9.035976297124169, 12.635875005705671, 14.418750173000687, 1.965008885765979, 0.15283054202871843, 4.502758034331163, 6.454043093352853, 7.963742350028338, 8.027052600914764, 10.967112294295905, -1.2430798571267325, 9.637785107184861, -4.185444591588373, 14.057581237090247, 4.205307322509485, 7.017243538340267, -2.728310162701156, 22.88486646364003, 62.445787169865454, 14.92846924337293, 100.99620624453172, -21.01775646758702, -11.15044046930089, 29.60692127146112, 79.1522449103986, 30.667134507660307, 33.2085825845372, 4.428035421125756, 8.34579839721005, -13.594798608401991, 11.870630835089448, 10.524744505153002, 12.85607548333148, 27.456035792398005, 22.058636447703805, 25.369478407887996, 31.83273816191705, 51.6239107853566, 96.67256936560901, 81.96836966182568, 38.000780842805796, 45.17462656350512, 81.02807589701786, -12.68394152521525, 81.59382369416979, 34.892024252939095, 21.381873093746858, 26.6898730927101, 43.88052032002462, 58.69283858569574, 54.49221359105495, 90.70663429834136, 8.057738525756491, -11.309313322303531, 149.9696988347689, 82.04954840310013, 70.69321799755184, 48.07921836861895, 51.25486691648695, 4.789788441271174, 76.83929269147256, 39.07468499703156, 16.654586750826763, 9.168238978750626, 13.28828219434635, 6.255172248358738, 2.3983693174398235, 0.6226764643217437, 5.010000035402371, -6.603243649643, 8.847824134689734, -2.876810983999812, 0.7289842837420277, -3.279758008755326, 5.220328211977897, 6.045818164917872, 0.8467414303055651, 10.281245726770265, 14.548947996130225, 20.65414895337922, 10.85002491779273, 48.14222919523364, 37.333154111143884, -11.923592311347413, 84.34912259516973, 36.48971880613473, 6.629687128247255, 102.3591405134866, 48.34000550201765, 101.42158253113894, 11.900038937739112, -55.164177107977885, 82.97776500627177, 42.62997530828501, 32.25884538887614, 10.089739839053118, 2.9171849752851178, -1.5054024057408553, 0.9810516970075019, 0.17783921474722353, 2.455688890263451, 2.5209650809638458, 43.67136869689038, 55.41437235759426, 22.932818553781573, 96.35868371238821, 9.747147943375012, -73.1014150786105, 90.85902103779652, 48.625433657605825, 45.742383021999984, 21.96689301306667, 17.178526830136214, 33.4189997729633, 17.83360055052242, -1.2012282113697452, 137.6529493657726, 78.18520718641314, 93.74203968862375, 137.18966810076998, 48.844121850945385, -6.68759587158339, 140.49399980921865, 61.815555091854705, 92.70675957628053, 35.99511904554866, 56.14541511011397, 39.81933510932956, 23.032971849978793, 38.484868807966265, 10.591570008600257, -21.053853002964225, 54.317726981909914, 17.736994007508446, 10.24222477852025, 14.353666201526924, 2.6522529282787515, 7.7936898510074935, 24.295435810543896, 26.765346672818794, 21.000931483121946, 50.10486183278114, -11.208611569072204, -10.404789846072106, 32.854852725233535, 34.53332998120255, 20.54591174808724, 15.048224903140508, -5.719938754180203, 8.648548614259598, 6.704276487626362, 11.747295787191081, 7.869879894559801, 18.43688919209728, 0.481985055793956, -14.613247780029923, 20.571616697026695, 8.738871127355821, 5.450928259405285, 4.454183167999255, -1.8926422365657665, 35.211512271180766, 16.134215741640418, 34.689377645757965, 17.31862858438992, 5.757139426675138, 7.936813467215693, 60.549662952987575, 59.97838400928519, 46.785473248181134, -5.759589951209215, 26.451987729513107, 11.585024235400546, 39.45981654852574, 54.55049304579518, 42.50165691823445, 53.50420620296736, 36.54001488731066, -14.457335163962142, 22.99693410426017, 15.874055394573483, 43.86531070314011, 18.998111027559478, 10.216708714382515, -2.1688154630695147, 9.30722706654009, -6.86186042810584, 4.194115990233447, 2.4104280459755283, 5.949323180192973, 1.7396956994814392, 12.342563462758726, -6.853514567646357, 5.794459899484764, -4.276905120697153, 16.7451891751795, -3.1887387927617126, 4.931128976121094, -4.047461569276654, 1.2256132895263046, 1.3845217294259675, -2.4414264424521335, 4.031243450889584, -5.647855401984925, 2.705870975747417, -0.927362293507838, -1.3579713026872344, 4.83524709778491, -3.706493353178872, 5.662867262929591, -3.8645897750975413, -4.277185248832327, 5.256281878860347, -1.3792776923857717, -1.5153871093967497, 20.43606059806084, 2.1423425907507117, 17.261314320142652, 14.343504929507482, 37.501092170756706, 43.42733246372812, 33.239205758811735, 110.4293122212254, 18.38992259318576, -73.81680759690288, 94.61169623650461, 38.88785958826219, 61.8586380765901, -1.7725944891724055, 13.985122702558229, -1.3917958749710468, 27.50634331359692, 7.073048414733424, 12.760674607101008, 7.938345130796581, 31.13575762275193, 3.694458029411912, 48.75417786313794, 1.757429436918855, -23.910377936456197, 21.844119324171768, 51.524738739852744, 38.016964063168246, 60.5616691931516, 26.27199047288113, 21.510686413296554, 17.154700921326608, 30.332925134201254, 27.342231232163606, 45.28177614923667, 22.116915225428237, 7.453141959589815, 8.114208572100477, -4.580341359907615, 17.783602229620474, -0.0567297233098536, 10.897381691153988, 19.12952856010387, -1.4002882585208738, -9.18870603703378, 11.730505697009248, 25.491581367620636, 9.229621381279722, 7.414159952282015, 7.048396913753095, -4.436328988635218, 0.9856379397344791, 10.23915378847725, 4.109183669289458, 42.016124163132936, 55.50854878870708, 40.50159586055089, 84.92050192589414, 47.722660272545, -39.108069068640106, 94.14874742850442, 97.12527322642819, 64.59723260509264, 77.92756983230991, 34.67077872547193, -35.00038564686949, 59.435459453702364, 27.810512561770757, 27.696261504423152, 24.597863891436063, 16.191523749616277, 16.434005652379096, 16.251228362291556, -1.4470092916392758, -6.807820981198831, 17.23313609527837, 27.17459045553477, 15.998366920801216, 34.04588720808663, 12.812853692097505, 76.65887309463474, 56.20645195523893, 37.79703582008142, 65.18998336355779, 51.75725712443415, 13.892195104008664, 55.33230090073317, 54.08173329842155, 52.54251732604313, 31.325533285346097, 15.430976072705857, 2.7669014323851293, -5.011245656345591, 5.944114762774313, 8.243096367334749, 16.400015588130746, 9.424160799237693, 19.043316574423116, -1.7987099195769254, -0.024711022455818465, 7.483486596346538, 16.52208029738839, 14.232722713373452, 7.3138183339509, 59.8760464756322, 33.640271577803006, 60.077260205678726, 30.47544040412973, 0.7202422361520462, -14.140411755277317, 75.69392795761885, 44.95043789730953, 44.9510239268608, 32.3205035262706, 23.495985940605244, 60.036427104529146, 54.353667128935214, 25.383878863265338, 57.91335526798447, 72.31073649004821, -26.738725512962368, 60.012759185821274, 30.53199492439031, 47.27512787065463, 17.58125379269688, 8.756351614737278, 1.4912993504313288, 0.7803865444325426, -0.6655762371834832, 1.0723978918996206, 1.9514052114174394, -4.1346568361112395, 1.4644368218658677, 5.2046842101369695, -4.708877608835773, 4.637302726898053, -1.924178738342968, 31.810130368930857, 6.025110331441391, 32.986430965573014, 38.46587157805113, 6.161598820644798, 29.65212624609925, 23.524518431697853, 55.56135708856271, 53.54322118139769, 41.823251424192044, 58.00829692319393, 8.089147826777037, -13.900205821468827, 44.64239652776174, 16.912290644739706, 27.582280942214744, -3.0102051484076, 11.294317205025395, -1.7591308424896237, 1.7981321618149249, 8.01233743660623, -0.63590773910407, 14.426603768792642, 50.46149870658108, 33.32747142684738, 65.78449779848714, 22.569761955026905, -27.478815134945677, 19.59036609855136, 52.68232341363446, 37.60924855344096, 15.604224395583074, 6.323819894089909, -0.5838114365476819, 3.649127557775844, -0.010999953185894529, 36.66437609743399, 35.821493987212676, 35.07250309828304, 59.02357386730795, 30.324682549027813, -24.093874606389768, 77.08519142776339, 83.23532917387968, 35.77818173832074, -4.4044568221111575, 65.22278970105597, 30.26153078877083, 42.04511536526688, 7.072949466758018, 5.842236035192645, 29.801036833362712, 15.459040643291559, 14.250334009892454, 6.671966880370303, 4.587141465515315, 2.2583001564250536, -3.0494366205834145, 2.617853080148029, 0.08520000539046829, 1.4863069775730438, 0.1806938021607768, -0.8211345639809078, 18.87088465365304, 10.557023181065958, 14.258362732580004, 18.28990230061076, 3.7631464977283855, -13.411136691002318, 23.865670351855684, 28.734359993938284, 7.230283170913205, 15.661314642081933, -5.4565950230765745, 3.7902361153238715, 2.134455091833335, 5.568144959135699, 7.775781904380718, 5.650357319179651, 0.3491747577830231, 9.241972980588896, -0.16136641727198509, 17.602254818503802, -4.778121832204973, 7.236277621500044, 4.409027202351213, 17.28656480022201, 8.512944204913904, 15.179408380090154, 9.627599877346988, 2.929818167969806, -1.4411498061786658, 9.442034793897768, 1.0697460660842402, 10.052409308165332, 2.983957907784283, 34.856708225466356, 19.62476938036685, 59.27053292467033, -6.567253900430611, 21.405415954084532, 4.44745498349544, 30.057805215141414, 39.81292085723175, 23.733933609887714, 10.569918625924306, 1.8316809843897943, 3.2872446737086225, 54.071443684311895, 37.38159239099255, 81.71289493018406, 34.41140676155446, 23.401323354367605, 35.74540562881364, 72.99093597361664, 74.95402903647124, 83.9308854129533, 13.05943965691904, 2.6345316824455547, 34.91205892292447, 20.273797603839455, 20.488818771635604, -4.441147509349168, 9.494357372154258, 2.2006070826679416, 8.378216493507466, 23.26830036196028, 30.219658264377646, 16.44840393465428, 60.01618023636568, 54.50371891683548, -2.9986077405345526, 71.72306682345092, 121.40479750737657, 36.08742688953561, -76.28773516854204, 80.76509188787686, 42.948841110906905, 29.12244419124259, 13.436010524266916, 1.9102232372163268, 4.795959222110396, -5.428597377293045, 8.100266115808273, -5.155666252433029, 13.563090895512044, 46.961347216638174, 91.55251465373894, 71.24503844934188, 140.72212893485195, -2.5620093628396887, -19.671503022459206, 52.066296272161, 108.78706336959814, 79.73198225234293, 31.607538422843312, 12.234366297528444, 2.6732829467381665, -0.46237146886054437, 1.0299726694240576, -3.0072730970852737, 0.9326611245549447, -1.850338615488982, 5.621350250061769, -5.569294657104442, 4.050659743575237, -4.0279088768440525, 3.7895522264974355, -3.110558663601048, 1.0252125628004185, 1.1519672177964284, 9.306880106675997, -0.16767630325076155, 9.85561016575597, -4.702095667961949, -3.617783892936926, 5.282710662150597, 10.271700769006289, -2.2203816554330817, 5.313322771491491, -2.969538158958421, 3.4729076882719823, -6.605655985437618, 4.14358337902802, -2.7006919211122726, 1.735138893753669, -1.635504574192461, 2.3480201016232924, -2.200654920293327, 3.6312841991766636, -1.7509421007308692, 5.8793746492191765, 51.21260644194571, 66.3202448438079, 68.18745516974442, 72.12451855615555, 12.958404063761918, -13.564260356179773, 63.47697907608831, 78.71962770209058, 56.9349094205497, 25.005680581937938, 2.4637648109153467, 0.8338327954572197, -0.2369548243519446, 2.915431542310241, 15.551584354753864, 13.660909299404675, 17.877018569037215, 6.039152151281462, -1.0138663197158166, 8.89649689963253, 34.25361175354402, 22.80253277212353, 36.42649845716956, 0.8675597653573561, -16.89426804732522, 19.126199233315884, 17.346073316071312, 45.68828969937242, 30.7061373740989, 29.092088801283523, 52.89987345525682, 1.5838602487329467, -12.005627283098221, 47.26009261119815, 32.93111964219945, 20.81680121828436, 16.87908046297324, 9.705051872401405, 71.05194439376906, 21.256349806211716, 104.97814443026395, 2.5878982039463025, 30.389067881249332, 27.06165256983793, 92.38117918550535, 59.95796018356863, 47.958324219172034, 56.00404181816687, 50.19383076895758, 37.69801727005461, 102.39033704073475, 25.844482968696326, -39.821926514568645, 79.69343114739122, 35.94611904096479, 8.133759321619952, 36.14179811430073, 22.03460517964289, 5.227757763263131, 5.384747369175224, 1.758214922780984, -3.681253236439444, 3.607183725846733, 2.1933911144916123, -4.805831108555243, 4.864828458158442, -2.035702446332937, 9.64273084295649, 23.687600361225435, 22.02408566505622, 24.592135307833583, 15.361180991960527, -13.01513141096757, 9.079447050807644, 35.818703487645465, 17.606840696503735, 17.69297385948232, 9.04755438997912, 10.402352213224912, 12.438477399777415, 22.130251380859747, 73.88402586718725, 2.252544788890175, 110.75122945907147, -29.88610917951408, 25.86019999156167, 40.52897194245768, 71.97964524376056, 49.65204248714802, 35.6079687369375, 5.65281685214282, -2.2938416004220343, 6.8877850599252195, 9.29374996687732, 11.820380583376465, 26.98552461532928, 16.144093202919024, 33.038662239126154, 8.543791025508995, -18.390953671648177, 21.3572781949627, 22.4201708292367, 20.88210494676044, 7.051286704151512, 1.3912163108822835, 3.693194263881428, 0.8930576817829116, 16.8618746801639, 28.657241772236212, 12.660625893559498, 63.88027619330255, 7.770328030331854, -19.63682180730467, 39.14311786517867, 32.15165027528068, 7.019652083407133, 25.480133286211753, 6.716786391382345, 35.214383063325286, 4.640307531623613, 32.5343440674102, 2.0681584317530426, 14.20823888689796, 18.8076663731222, 38.00235832520856, 20.59614377698689, 1.617003987746763, 17.406871197443646, 29.86234325326896, 23.033399772766003, 19.324091357745907, 15.472303054468178, 1.3753164240261668, 6.5558679314167465, 10.872007224715496, 28.404053359862]"
$ws.Range("E2").Value = 20.00055902852358
$ws.Range("F2").Value = 28.69137933855449

$ws.Rows.Item(2).AutoFit()
